# Data refresh for the FFXIV Leve-profit tracking sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Each worksheet is a Table_<JOB> with columns:
#   H currentAveragePrice   I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ           L LevePriceHQ             M LeveProfitNQ   N LeveProfitHQ
# The scheduled runner repriced a batch of leve rows; update the affected cells in place.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1550.6666
$ws.Range("I29").Value = 64.5
$ws.Range("K29").Value = 193.5
$ws.Range("M29").Value = 87.5
$ws.Range("H38").Value = 1684
$ws.Range("I38").Value = 69.42856999999999
$ws.Range("J38").Value = 2711.4546
$ws.Range("K38").Value = 208.28571
$ws.Range("L38").Value = 8134.3638
$ws.Range("M38").Value = 163.71429
$ws.Range("N38").Value = -8878.363799999999
$ws.Range("H58").Value = 1310.8235
$ws.Range("J58").Value = 2211.111
$ws.Range("L58").Value = 6633.333
$ws.Range("N58").Value = -6933.333
$ws.Range("H87").Value = 24766.666
$ws.Range("J87").Value = 24766.666
$ws.Range("L87").Value = 24766.666
$ws.Range("N87").Value = -27262.666
$ws.Range("H90").Value = 24766.666
$ws.Range("J90").Value = 24766.666
$ws.Range("L90").Value = 74299.99800000001
$ws.Range("N90").Value = -86779.99800000001
$ws.Range("H135").Value = 1269.5714
$ws.Range("I135").Value = 1102.6111
$ws.Range("J135").Value = 1570.1
$ws.Range("K135").Value = 9923.499900000001
$ws.Range("L135").Value = 14130.9
$ws.Range("M135").Value = -7388.499900000001
$ws.Range("N135").Value = -19200.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4580.125
$ws.Range("I97").Value = 5091.5713
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 5091.5713
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -4595.5713
$ws.Range("N97").Value = -1992
$ws.Range("H110").Value = 1194.2858
$ws.Range("I110").Value = 817.05884
$ws.Range("J110").Value = 2797.5
$ws.Range("K110").Value = 817.05884
$ws.Range("L110").Value = 2797.5
$ws.Range("M110").Value = 1227.94116
$ws.Range("N110").Value = -6887.5
$ws.Range("H115").Value = 27650
$ws.Range("J115").Value = 27650
$ws.Range("L115").Value = 27650
$ws.Range("N115").Value = -30784

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12713.588
$ws.Range("I94").Value = 906.3570999999999
$ws.Range("J94").Value = 67814
$ws.Range("K94").Value = 906.3570999999999
$ws.Range("L94").Value = 67814
$ws.Range("M94").Value = -455.3570999999999
$ws.Range("N94").Value = -68716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2741657.2
$ws.Range("I31").Value = 4546799.5
$ws.Range("J31").Value = 2820.8276
$ws.Range("K31").Value = 4546799.5
$ws.Range("L31").Value = 2820.8276
$ws.Range("M31").Value = -4546504.5
$ws.Range("N31").Value = -3410.8276
$ws.Range("H34").Value = 2741657.2
$ws.Range("I34").Value = 4546799.5
$ws.Range("J34").Value = 2820.8276
$ws.Range("K34").Value = 4546799.5
$ws.Range("L34").Value = 2820.8276
$ws.Range("M34").Value = -4546597.5
$ws.Range("N34").Value = -3224.8276
$ws.Range("H70").Value = 30090
$ws.Range("J70").Value = 30090
$ws.Range("L70").Value = 30090
$ws.Range("N70").Value = -30720
$ws.Range("H73").Value = 30090
$ws.Range("J73").Value = 30090
$ws.Range("L73").Value = 30090
$ws.Range("N73").Value = -32274
$ws.Range("H81").Value = 31031.5
$ws.Range("I81").Value = 13298
$ws.Range("J81").Value = 36942.668
$ws.Range("K81").Value = 13298
$ws.Range("L81").Value = 36942.668
$ws.Range("M81").Value = -12300
$ws.Range("N81").Value = -38938.668
$ws.Range("H84").Value = 31031.5
$ws.Range("I84").Value = 13298
$ws.Range("J84").Value = 36942.668
$ws.Range("K84").Value = 39894
$ws.Range("L84").Value = 110828.004
$ws.Range("M84").Value = -34902
$ws.Range("N84").Value = -120812.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2851.45
$ws.Range("I60").Value = 173
$ws.Range("J60").Value = 4293.6924
$ws.Range("K60").Value = 519
$ws.Range("L60").Value = 12881.0772
$ws.Range("M60").Value = -268
$ws.Range("N60").Value = -13383.0772
$ws.Range("H94").Value = 4165.875
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4165.875
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 12497.625
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -13849.625
$ws.Range("H124").Value = 2434.8
$ws.Range("I124").Value = 2615
$ws.Range("J124").Value = 2314.6667
$ws.Range("K124").Value = 7845
$ws.Range("L124").Value = 6944.000100000001
$ws.Range("M124").Value = -2935
$ws.Range("N124").Value = -16764.0001
$ws.Range("H141").Value = 793
$ws.Range("I141").Value = 703.3333
$ws.Range("J141").Value = 1600
$ws.Range("K141").Value = 2109.9999
$ws.Range("L141").Value = 4800
$ws.Range("M141").Value = 3070.0001
$ws.Range("N141").Value = -15160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 8800
$ws.Range("I15").Value = 7000
$ws.Range("K15").Value = 7000
$ws.Range("M15").Value = -6712
$ws.Range("H32").Value = 25000
$ws.Range("J32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("N32").Value = -25592
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H80").Value = 2305
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 2206.25
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 2206.25
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4202.25
$ws.Range("H81").Value = 8800
$ws.Range("I81").Value = 7000
$ws.Range("K81").Value = 7000
$ws.Range("M81").Value = -6002
$ws.Range("H83").Value = 2305
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 2206.25
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 11031.25
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -21015.25
$ws.Range("H84").Value = 8800
$ws.Range("I84").Value = 7000
$ws.Range("K84").Value = 21000
$ws.Range("M84").Value = -16008
$ws.Range("H102").Value = 2112
$ws.Range("I102").Value = 1963.5
$ws.Range("J102").Value = 3300
$ws.Range("K102").Value = 1963.5
$ws.Range("L102").Value = 3300
$ws.Range("M102").Value = -341.5
$ws.Range("N102").Value = -6544
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20452
$ws.Range("H55").Value = 280.30768
$ws.Range("I55").Value = 292.75
$ws.Range("J55").Value = 260.4
$ws.Range("K55").Value = 292.75
$ws.Range("L55").Value = 260.4
$ws.Range("M55").Value = -119.75
$ws.Range("N55").Value = -606.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 26033.334
$ws.Range("J27").Value = 26033.334
$ws.Range("L27").Value = 26033.334
$ws.Range("N27").Value = -26171.334
$ws.Range("H54").Value = 17999.143
$ws.Range("J54").Value = 17999.143
$ws.Range("L54").Value = 17999.143
$ws.Range("N54").Value = -19039.143
$ws.Range("H81").Value = 38463236
$ws.Range("I81").Value = 55557070
$ws.Range("J81").Value = 2101.25
$ws.Range("K81").Value = 111114140
$ws.Range("L81").Value = 4202.5
$ws.Range("M81").Value = -111113079
$ws.Range("N81").Value = -6324.5
$ws.Range("H84").Value = 38463236
$ws.Range("I84").Value = 55557070
$ws.Range("J84").Value = 2101.25
$ws.Range("K84").Value = 555570700
$ws.Range("L84").Value = 21012.5
$ws.Range("M84").Value = -555565396
$ws.Range("N84").Value = -31620.5
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
